$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 90.5
$ws.Range("I8").Value = 17.222221
$ws.Range("J8").Value = 750
$ws.Range("K8").Value = 51.666663
$ws.Range("L8").Value = 2250
$ws.Range("M8").Value = 87.333337
$ws.Range("N8").Value = -2528

$ws.Range("H58").Value = 2599.75
$ws.Range("J58").Value = 9333
$ws.Range("L58").Value = 27999
$ws.Range("N58").Value = -28299

$ws.Range("H62").Value = 73089.37
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 73089.37
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws.Range("H86").Value = 147650.86
$ws.Range("I86").Value = 337297.34
$ws.Range("K86").Value = 337297.34
$ws.Range("M86").Value = -336174.34

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").Value = $null

$ws.Range("H89").Value = 147650.86
$ws.Range("I89").Value = 337297.34
$ws.Range("K89").Value = 1686486.7
$ws.Range("M89").Value = -1680870.7

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").Value = $null

$ws.Range("H98").Value = 2264.4
$ws.Range("I98").Value = 1524.3125
$ws.Range("K98").Value = 1524.3125
$ws.Range("M98").Value = -26.3125

$ws.Range("H122").Value = 2264.4
$ws.Range("I122").Value = 1524.3125
$ws.Range("K122").Value = 4572.9375
$ws.Range("M122").Value = -2122.9375

$ws.Range("H137").Value = 3493.2354
$ws.Range("I137").Value = 2923.7144
$ws.Range("J137").Value = 3891.9
$ws.Range("K137").Value = 8771.143199999999
$ws.Range("L137").Value = 11675.7
$ws.Range("M137").Value = -6221.143199999999
$ws.Range("N137").Value = -16775.7

$ws.Range("H138").Value = 7289.754
$ws.Range("I138").Value = 9080
$ws.Range("J138").Value = 6842.1924
$ws.Range("K138").Value = 27240
$ws.Range("L138").Value = 20526.5772
$ws.Range("M138").Value = -22100
$ws.Range("N138").Value = -30806.5772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 9936.166999999999
$ws.Range("I16").Value = 179.5
$ws.Range("K16").Value = 179.5
$ws.Range("M16").Value = 107.5

$ws.Range("H32").Value = 34370.5
$ws.Range("I32").Value = 21935.309
$ws.Range("J32").Value = 52332.445
$ws.Range("K32").Value = 21935.309
$ws.Range("L32").Value = 52332.445
$ws.Range("M32").Value = -21648.309
$ws.Range("N32").Value = -52906.445

$ws.Range("H45").Value = 1252586.1
$ws.Range("I45").Value = 2001038.2
$ws.Range("K45").Value = 2001038.2
$ws.Range("M45").Value = -2000661.2

$ws.Range("H61").Value = 8277.546
$ws.Range("I61").Value = 8105.3
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 8105.3
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -7893.3
$ws.Range("N61").Value = -10424

$ws.Range("H80").Value = 61419.4
$ws.Range("J80").Value = 56774.5
$ws.Range("L80").Value = 56774.5
$ws.Range("N80").Value = -58770.5

$ws.Range("H83").Value = 61419.4
$ws.Range("J83").Value = 56774.5
$ws.Range("L83").Value = 170323.5
$ws.Range("N83").Value = -180307.5

$ws.Range("H132").Value = 29328.21
$ws.Range("I132").Value = 31364.742
$ws.Range("J132").Value = 5568.6665
$ws.Range("K132").Value = 94094.226
$ws.Range("L132").Value = 16705.9995
$ws.Range("M132").Value = -91564.226
$ws.Range("N132").Value = -21765.9995

$ws.Range("H136").Value = 8277.546
$ws.Range("I136").Value = 8105.3
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 24315.9
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -21765.9
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2876.077
$ws.Range("I105").Value = 2199.3333
$ws.Range("J105").Value = 3456.1428
$ws.Range("K105").Value = 2199.3333
$ws.Range("L105").Value = 3456.1428
$ws.Range("M105").Value = -452.3332999999998
$ws.Range("N105").Value = -6950.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 746.8
$ws.Range("I94").Value = 624
$ws.Range("J94").Value = 828.6667
$ws.Range("K94").Value = 624
$ws.Range("L94").Value = 828.6667
$ws.Range("M94").Value = -173
$ws.Range("N94").Value = -1730.6667

$ws.Range("H118").Value = 97996.5
$ws.Range("J118").Value = 97996.5
$ws.Range("L118").Value = 97996.5
$ws.Range("N118").Value = -101310.5

$ws.Range("H122").Value = 1899
$ws.Range("I122").Value = 1899
$ws.Range("K122").Value = 5697
$ws.Range("M122").Value = -3247

$ws.Range("H131").Value = 39298.5
$ws.Range("J131").Value = 39099.2
$ws.Range("L131").Value = 39099.2
$ws.Range("N131").Value = -49179.2

$ws.Range("H141").Value = 526313.1
$ws.Range("J141").Value = 549908.5600000001
$ws.Range("L141").Value = 549908.5600000001
$ws.Range("N141").Value = -560268.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1372.6428
$ws.Range("I68").Value = 1107.25
$ws.Range("J68").Value = 1726.5
$ws.Range("K68").Value = 3321.75
$ws.Range("L68").Value = 5179.5
$ws.Range("M68").Value = -2510.75
$ws.Range("N68").Value = -6801.5

$ws.Range("H71").Value = 1372.6428
$ws.Range("I71").Value = 1107.25
$ws.Range("J71").Value = 1726.5
$ws.Range("K71").Value = 9965.25
$ws.Range("L71").Value = 15538.5
$ws.Range("M71").Value = -5909.25
$ws.Range("N71").Value = -23650.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 15713.286
$ws.Range("J21").Value = 15999
$ws.Range("L21").Value = 15999
$ws.Range("N21").Value = -16345

$ws.Range("H30").Value = 15713.286
$ws.Range("J30").Value = 15999
$ws.Range("L30").Value = 15999
$ws.Range("N30").Value = -16209

$ws.Range("H80").Value = 3499.5
$ws.Range("I80").Value = 3499
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 3499
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -2501
$ws.Range("N80").Value = -5496

$ws.Range("H83").Value = 3499.5
$ws.Range("I83").Value = 3499
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 17495
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -12503
$ws.Range("N83").Value = -27484

$ws.Range("H113").Value = 52695.15
$ws.Range("I113").Value = 69421.2
$ws.Range("J113").Value = 2517
$ws.Range("K113").Value = 69421.2
$ws.Range("L113").Value = 2517
$ws.Range("M113").Value = -67251.2
$ws.Range("N113").Value = -6857

$ws.Range("H122").Value = 3590.5
$ws.Range("I122").Value = 3308.6
$ws.Range("K122").Value = 9925.799999999999
$ws.Range("M122").Value = -7475.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 3750
$ws.Range("I10").Value = 3750
$ws.Range("K10").Value = 3750
$ws.Range("M10").Value = -3610

$ws.Range("H40").Value = 4046.125
$ws.Range("I40").Value = 3158.5454
$ws.Range("K40").Value = 3158.5454
$ws.Range("M40").Value = -3022.5454

$ws.Range("H61").Value = 3979.8823
$ws.Range("I61").Value = 3706.5386
$ws.Range("K61").Value = 3706.5386
$ws.Range("M61").Value = -3504.5386

$ws.Range("H113").Value = 3979.8823
$ws.Range("I113").Value = 3706.5386
$ws.Range("K113").Value = 3706.5386
$ws.Range("M113").Value = -1536.5386

$ws.Range("H132").Value = 32209.977
$ws.Range("I132").Value = 36762.36
$ws.Range("J132").Value = 4895.6665
$ws.Range("K132").Value = 110287.08
$ws.Range("L132").Value = 14686.9995
$ws.Range("M132").Value = -107757.08
$ws.Range("N132").Value = -19746.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 590130.5
$ws.Range("J81").Value = 3334550
$ws.Range("L81").Value = 6669100
$ws.Range("N81").Value = -6671222

$ws.Range("H84").Value = 590130.5
$ws.Range("J84").Value = 3334550
$ws.Range("L84").Value = 33345500
$ws.Range("N84").Value = -33356108

$ws.Range("H107").Value = 590.5238000000001
$ws.Range("I107").Value = 440.22223
$ws.Range("K107").Value = 1320.66669
$ws.Range("M107").Value = 599.33331

$ws.Range("H116").Value = 89998
$ws.Range("J116").Value = 89998
$ws.Range("L116").Value = 89998
$ws.Range("N116").Value = -99176

$ws.Range("H122").Value = 3336.353
$ws.Range("I122").Value = 2741.2
$ws.Range("J122").Value = 7800
$ws.Range("K122").Value = 8223.599999999999
$ws.Range("L122").Value = 23400
$ws.Range("M122").Value = -5773.599999999999
$ws.Range("N122").Value = -28300

$ws.Range("H126").Value = 32602.688
$ws.Range("I126").Value = 53869.26
$ws.Range("J126").Value = 1520.7693
$ws.Range("K126").Value = 161607.78
$ws.Range("L126").Value = 4562.3079
$ws.Range("M126").Value = -159137.78
$ws.Range("N126").Value = -9502.3079

$ws.Range("H132").Value = 38138.555
$ws.Range("I132").Value = 40416.133
$ws.Range("J132").Value = 26750.666
$ws.Range("K132").Value = 121248.399
$ws.Range("L132").Value = 80251.99800000001
$ws.Range("M132").Value = -118718.399
$ws.Range("N132").Value = -85311.99800000001

$ws.Range("H136").Value = 34306736
$ws.Range("I136").Value = 4074480
$ws.Range("K136").Value = 12223440
$ws.Range("M136").Value = -12220890
